$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "locial" -> "logical" in the Purpose row (A5).
# The engine will automatically retire the old shared-string entry and
# append the corrected text as a new shared-string at the end - matching
# the target sharedStrings.xml layout.
$ws.Range("A5").Value() = "Purpose: Unit test the logical structure of the Check_State Class and its Interface"

# Move the view: scroll back so column A is visible again (clears any
# topLeftCell scroll offset) and select A7:E7 (the Comments row).
$ws.Range("A7:E7").Select()
